$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename category "Family Policy" -> "Parental Leave Reform" for rows 26-30 ---
$ws.Range("D26").Value = "Parental Leave Reform"
$ws.Range("D27").Value = "Parental Leave Reform"
$ws.Range("D28").Value = "Parental Leave Reform"
$ws.Range("D29").Value = "Parental Leave Reform"
$ws.Range("D30").Value = "Parental Leave Reform"

# --- Row 31 (coronavirusLockdownR1): clear the (wrong-sign / unreliable) average age value,
#     rewrite the notes explaining why, and grow the row to fit the longer note ---
$ws.Range("E31").ClearContents()
$ws.Range("H31").Value = "Hard to tell what the average age of beneficiaries is. Age of people who die of covid19 is 81 according to RKI. But most of the WTP is due to lower economic costs of a smaller R. Missing Value for now."
$ws.Rows.Item(31).RowHeight = 120

# --- Move the view / selection from the previous edit location to the new one ---
$ws.Activate()
$ws.Range("D21").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
